$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Size 3")

# Update contribution values to 1 (marking as "Agreed") for the listed cells.
# The G column formula recalculates automatically to "Ok" once the row sums to 1.
$ws.Range("D5").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("D24").Value = 1

# Update the active selection to match the final cursor location.
$ws.Activate()
$ws.Range("D25").Select()
